$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D19").Value = "Число парков культуры - parks (шт.) (8017015)"
$ws.Range("D20").Value = "Число театров - theatres (шт.) (8017007)"

# Match the style (fill/border/center alignment) already used by D18
$ws.Range("D18").Copy()
$ws.Range("D19:D20").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("E25").Select()
